$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format temporarily so numeric-looking strings
# (e.g. "584.50", "1.00") are stored as text, not converted to numbers.
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Range("D2").Value = '66.522.49'
$ws.Range("E2").Value = '  -0.86%  '
$ws.Range("D3").Value = '2.576.74'
$ws.Range("E3").Value = '  -2.02%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '584.50'
$ws.Range("E5").Value = '  -1.63%  '
$ws.Range("D6").Value = '167.84'
$ws.Range("E6").Value = '  +0.77%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '0.527'
$ws.Range("E8").Value = '  -1.22%  '
$ws.Range("D9").Value = '2.575.77'
$ws.Range("E9").Value = '  -2.08%  '
$ws.Range("D10").Value = '0.139'
$ws.Range("E10").Value = '  -1.12%  '
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("D12").Value = '0.355'
$ws.Range("E12").Value = '  -1.25%  '
$ws.Range("D13").Value = '5.16'
$ws.Range("E13").Value = '  -1.29%  '
$ws.Range("D14").Value = '26.75'
$ws.Range("E14").Value = '  -3.55%  '
$ws.Range("E15").Value = '  -2.31%  '
$ws.Range("D16").Value = '0.0000177'
$ws.Range("E16").Value = '  -2.19%  '
$ws.Range("D17").Value = '66.354.09'
$ws.Range("E17").Value = '  -1.02%  '
$ws.Range("D18").Value = '2.565.50'
$ws.Range("E18").Value = '  -2.30%  '
$ws.Range("D19").Value = '11.40'
$ws.Range("E19").Value = '  -6.37%  '
$ws.Range("D20").Value = '7.76'
$ws.Range("E20").Value = '  -3.82%  '
$ws.Range("D21").Value = '350.90'
$ws.Range("E21").Value = '  -1.80%  '
$ws.Range("D22").Value = '4.23'
$ws.Range("E22").Value = '  -1.98%  '
$ws.Range("D23").Value = '4.61'
$ws.Range("E23").Value = '  -1.34%  '
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("D25").Value = '1.90'
$ws.Range("E25").Value = '  -1.20%  '
$ws.Range("D26").Value = '69.12'
$ws.Range("E26").Value = '  -1.59%  '
$ws.Range("D27").Value = '9.85'
$ws.Range("E27").Value = '  -9.79%  '
$ws.Range("E28").Value = '  -1.62%  '
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D30").Value = '0.0₃0988'
$ws.Range("E30").Value = '  -2.22%  '
$ws.Range("D31").Value = '528.10'
$ws.Range("E31").Value = '  -4.03%  '
$ws.Range("D32").Value = '8.14'
$ws.Range("E32").Value = '  +2.79%  '
$ws.Range("D33").Value = '1.33'
$ws.Range("E33").Value = '  -2.47%  '
$ws.Range("D34").Value = '1.85'
$ws.Range("E34").Value = '  -2.95%  '
$ws.Range("E35").Value = '  -3.35%  '
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  -0.09%  '
$ws.Range("D37").Value = '1.46'
$ws.Range("E37").Value = '  -3.08%  '
$ws.Range("D38").Value = '156.09'
$ws.Range("E38").Value = '  +0.39%  '
$ws.Range("D39").Value = '18.77'
$ws.Range("E39").Value = '  -1.81%  '
$ws.Range("D40").Value = '0.359'
$ws.Range("E40").Value = '  -1.90%  '
$ws.Range("D41").Value = '18.32'
$ws.Range("E41").Value = '  +2.14%  '
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("D45").Value = '2.43'
$ws.Range("E45").Value = '  -0.05%  '
$ws.Range("E46").Value = '  -3.83%  '
$ws.Range("D47").Value = '149.17'
$ws.Range("E47").Value = '  -1.60%  '
$ws.Range("D48").Value = '0.567'
$ws.Range("E48").Value = '  -2.51%  '
$ws.Range("D49").Value = '3.72'
$ws.Range("E49").Value = '  -1.73%  '
$ws.Range("D50").Value = '1.72'
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("E51").Value = '  -1.18%  '

# Swap RenderToken (was row 42) and Stacks (was row 43), with updated values
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '1.78'
$ws.Range("E42").Value = '  -0.61%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D43").Value = '5.11'
$ws.Range("E43").Value = '  -1.05%  '

# Restore original (default/general) style for column D so no residual
# text-format styling is left behind on the cells.
$colD.Style = "Normal"
